$d = $word.ActiveDocument

# --- Helper: replace a paragraph's content, keeping its paragraph mark ---
function Set-ParaText($para, [string]$text) {
    $r = $d.Range($para.Range.Start, $para.Range.End - 1)
    $r.Text = $text
}

# --- Helper: insert a brand-new paragraph right after $afterPara and return it ---
function New-ParaAfter($afterPara) {
    $afterPara.Range.InsertParagraphAfter()
    return $afterPara.Next()
}

# IMPORTANT: this engine resolves a held Paragraph reference by its *current*
# ordinal position, not by stable identity. So we never pre-fetch paragraphs
# far ahead via Paragraphs.Item(N) and hold onto them across intervening
# inserts elsewhere in the doc - instead we always walk strictly left to
# right with .Next()/New-ParaAfter(), fetching 'the next original paragraph'
# only once everything earlier has already been edited/inserted.

# =================================================================
# Section 1: 'EGD Findings' -> 'Indications' (+ Medications / Monitoring)
# =================================================================
$pHeading1 = $d.Paragraphs.Item(2)   # 'EGD Findings' -> 'Indications'
$pBody1    = $pHeading1.Next()        # EGD narrative -> Indications body
Set-ParaText $pHeading1 'Indications'
Set-ParaText $pBody1 '69-year-old male patient is here for an ERCP procedure for therapy of a pancreatic duct stricture and stent.'

$pMedH = New-ParaAfter $pBody1
$pMedH.Style = "Heading3"
Set-ParaText $pMedH 'Medications'

$pMedB = New-ParaAfter $pMedH
$pMedB.Style = "Normal"
Set-ParaText $pMedB 'Refer to record of source.'

$pMonH = New-ParaAfter $pMedB
$pMonH.Style = "Heading3"
Set-ParaText $pMonH 'Monitoring'

$pMonB = New-ParaAfter $pMonH
$pMonB.Style = "Normal"
Set-ParaText $pMonB 'Johns Hopkins Standard.'

# =================================================================
# Section 2: 'ERCP Findings' -> 'History' (+ Description of Procedure,
# Findings, ERCP Quality Metrics)
# =================================================================
$pHeading2 = $pMonB.Next()    # originally 'ERCP Findings' -> 'History'
$pBody2    = $pHeading2.Next() # ERCP narrative -> History body
Set-ParaText $pHeading2 'History'
Set-ParaText $pBody2 'The patient reports  history of therapy of a pancreatic duct stricture and stent.Patient presents with decreased appetite and abdominal pain, but no fever or leukocytosis. No history of no fever or leukocytosis is reported unless otherwise specified.Metastatic hepatosplenic carcinoma, chronic constipation, pancreatic duct stricture, stone, extensive tobacco use, recent acute pancreatitis.  extensive tobacco use, approximately 100 pack years.'

$pDescH = New-ParaAfter $pBody2
$pDescH.Style = "Heading2"
Set-ParaText $pDescH 'Description of Procedure'

$pDescB = New-ParaAfter $pDescH
$pDescB.Style = "Normal"
Set-ParaText $pDescB 'After the risks, benefits and alternatives of the procedure were thoroughly explained, informed consent was verified, confirmed and timeout was successfully executed by the treatment team. With the patient in the left semi-prone position, medications were administered intravenously. The duodenoscope Olympus TJF Q180B was passed from the mouth into the esophagus and further advanced from the esophagus into the stomach. From the stomach, the scope was directed to the second portion of the duodenum.'

$pFindH = New-ParaAfter $pDescB
$pFindH.Style = "Heading2"
Set-ParaText $pFindH 'Findings'

$pFindB = New-ParaAfter $pFindH
$pFindB.Style = "Normal"
Set-ParaText $pFindB 'A scout film of the abdomen was performed. It appeared normal.The duodenoscope Olympus TJF Q180B was advanced to the second portion of the duodenum without difficulty and without detailed examination of the upper GI tract. The esophagus, stomach, and duodenum appeared unremarkable on limited inspection.The major papilla was identified. It appeared normal in position and morphology, with an intact orifice and no surrounding erythema or edema.The minor papilla was identified and inspected. It appeared normal in position and morphology, with an intact orifice and no surrounding erythema or edema.The ampulla was identified and found to be normal.Bile duct cannulation was attempted using a sphincterotome preloaded with a guidewire (Visiglide 0.025 inch).Bile duct cannulation was unsuccessful. The procedure was terminated.Pancreatic duct cannulation was attempted. Pancreatic duct was selectively cannulated with standard wire-guided.Contrast was injected under fluoroscopic guidance and cholangiogram was performed. Multiple calcified calculi in the dorsal pancreatic duct.Pancreatogram was obtained and dilation, stones, incomplete pancreas divisum identified. Side branches were dilated. There was no communication with pseudocyst.Sphincterotomy was performed.Estimated blood loss: None.Specimens removed: None.Complications: There were no immediate complications.'

$pQualH = New-ParaAfter $pFindB
$pQualH.Style = "Heading2"
Set-ParaText $pQualH 'ERCP Quality Metrics'

$pQualB = New-ParaAfter $pQualH
$pQualB.Style = "Normal"
Set-ParaText $pQualB 'Difficulty of ERCP:  pancreatic duct cannulation successful.Cannulation success: Yes.Post-ERCP pancreatitis prophylaxis:Lactated ringers: Yes, refer to record of source.Rectal indomethacin:  Yes.Successful completion of intended procedure: Yes.Failed ERCP from another facility or provider: No.'

# =================================================================
# Section 3: Impressions list items 1-5 rewritten; 6-7 appended
# (the 'Impressions' heading itself is untouched)
# =================================================================
$pImpHeading = $pQualB.Next()  # 'Impressions' - unchanged
$pImp1 = $pImpHeading.Next()
$pImp2 = $pImp1.Next()
$pImp3 = $pImp2.Next()
$pImp4 = $pImp3.Next()
$pImp5 = $pImp4.Next()
Set-ParaText $pImp1 '1. Normal esophagus and stomach on limited views'
Set-ParaText $pImp2 '2. Normal major and minor papilla'
Set-ParaText $pImp3 '3. Successful pancreatic duct cannulation with sphincterotome and Visiglide wire'
Set-ParaText $pImp4 '4. Mildly dilated dorsal pancreatic duct with minimal side branch dilation'
Set-ParaText $pImp5 '5. Multiple calcified calculi in the dorsal pancreatic duct'

$pImp6 = New-ParaAfter $pImp5
Set-ParaText $pImp6 '6. Electrohydraulic lithotripsy performed to fragment pancreatic duct stone'

$pImp7 = New-ParaAfter $pImp6
Set-ParaText $pImp7 '7. Clear pancreatic duct with no remaining stones after procedure'

# =================================================================
# Section 4: New 'Recommendations' heading + 3 list items at the end
# =================================================================
$pRecH = New-ParaAfter $pImp7
$pRecH.Style = "Heading2"
Set-ParaText $pRecH 'Recommendations'

$pRec1 = New-ParaAfter $pRecH
$pRec1.Style = "Normal"
$pRec1.Format.SpaceAfter = 0
Set-ParaText $pRec1 '1. Finish IV fluids now.'

$pRec2 = New-ParaAfter $pRec1
Set-ParaText $pRec2 '2. Pain control as needed.'

$pRec3 = New-ParaAfter $pRec2
Set-ParaText $pRec3 '3. Follow up with referring provider.'

Write-Output "Edit complete."
